# Applies the MicroRTS Results.xlsx edit:
# - Adds two new worksheets "RawResults#5" and "BestResults#5" at the end,
#   populated with the new lap/results data.
# - Updates the active sheet / selection state so BestResults#5 ends up
#   the active tab, matching the recorded author session.

$wb = $excel.ActiveWorkbook

# --- Leave a new selection on RawResults#4 (mirrors the author's last click
#     on that sheet before adding the new ones) -----------------------------
$rawResults4 = $wb.Worksheets.Item("RawResults#4")
$rawResults4.Range("I9").Select() | Out-Null

# --- Add the two new worksheets at the end of the workbook -----------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$rawResults5 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$rawResults5.Name = "RawResults#5"

$bestResults5 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $rawResults5)
$bestResults5.Name = "BestResults#5"

# --- Populate RawResults#5 (single column A) --------------------------------
$rawData = @(
    ,@(2, 'A', 'Killer has cost: 4 hp: 3 min Damage: 1 max Damage: 1 attack range: 1')
    ,@(3, 'A', 'The number of games where killer was made was: 1')
    ,@(4, 'A', 'The number of games won by Player 0 with killer: 0')
    ,@(5, 'A', 'The score was: 0.05158371')
    ,@(6, 'A', 'Killer has cost: 2 hp: 3 min Damage: 1 max Damage: 1 attack range: 1')
    ,@(7, 'A', 'The number of games where killer was made was: 2')
    ,@(8, 'A', 'The number of games won by Player 0 with killer: 0')
    ,@(9, 'A', 'The score was: 0.06735747')
    ,@(10, 'A', 'Killer has cost: 3 hp: 4 min Damage: 1 max Damage: 1 attack range: 1')
    ,@(11, 'A', 'The number of games where killer was made was: 2')
    ,@(12, 'A', 'The number of games won by Player 0 with killer: 0')
    ,@(13, 'A', 'The score was: 0.02704026')
    ,@(14, 'A', 'Killer has cost: 3 hp: 2 min Damage: 1 max Damage: 1 attack range: 1')
    ,@(15, 'A', 'The number of games where killer was made was: 0')
    ,@(16, 'A', 'The number of games won by Player 0 with killer: 0')
    ,@(17, 'A', 'The score was: 0.0')
    ,@(18, 'A', 'Killer has cost: 3 hp: 3 min Damage: 2 max Damage: 2 attack range: 1')
    ,@(19, 'A', 'The number of games where killer was made was: 0')
    ,@(20, 'A', 'The number of games won by Player 0 with killer: 0')
    ,@(21, 'A', 'The score was: 0.0')
    ,@(22, 'A', 'Killer has cost: 3 hp: 3 min Damage: 1 max Damage: 2 attack range: 1')
    ,@(23, 'A', 'The number of games where killer was made was: 0')
    ,@(24, 'A', 'The number of games won by Player 0 with killer: 0')
    ,@(25, 'A', 'The score was: 0.0')
    ,@(26, 'A', 'Killer has cost: 3 hp: 3 min Damage: 0 max Damage: 1 attack range: 1')
    ,@(27, 'A', 'The number of games where killer was made was: 0')
    ,@(28, 'A', 'The number of games won by Player 0 with killer: 0')
    ,@(29, 'A', 'The score was: 0.0')
    ,@(30, 'A', 'Killer has cost: 3 hp: 3 min Damage: 0 max Damage: 0 attack range: 1')
    ,@(31, 'A', 'The number of games where killer was made was: 0')
    ,@(32, 'A', 'The number of games won by Player 0 with killer: 0')
    ,@(33, 'A', 'The score was: 0.0')
    ,@(34, 'A', 'Killer has cost: 3 hp: 3 min Damage: 1 max Damage: 1 attack range: 2')
    ,@(35, 'A', 'The number of games where killer was made was: 0')
    ,@(36, 'A', 'The number of games won by Player 0 with killer: 0')
    ,@(37, 'A', 'The score was: 0.0')
    ,@(38, 'A', 'Killer has cost: 3 hp: 3 min Damage: 1 max Damage: 1 attack range: 1')
    ,@(39, 'A', 'The number of games where killer was made was: 1')
    ,@(40, 'A', 'The number of games won by Player 0 with killer: 0')
    ,@(41, 'A', 'The score was: 0.016931217')
    ,@(42, 'A', 'Killer has cost: 4 hp: 3 min Damage: 1 max Damage: 1 attack range: 1')
    ,@(43, 'A', 'The number of games where killer was made by player 0 was: 1')
    ,@(44, 'A', 'The number of games where killer was made by player 1 was: 1')
    ,@(45, 'A', 'The number of games where killer was made by both players was: 0')
    ,@(46, 'A', 'The number of games won by Player 0 with killer: 0')
    ,@(47, 'A', 'The number of games won by Player 1 with killer: 2')
    ,@(48, 'A', 'The score was: 0.5')
)

foreach ($item in $rawData) {
    $row = $item[0]
    $col = $item[1]
    $val = $item[2]
    $rawResults5.Range("$col$row").Value = $val
}

# --- Populate BestResults#5 (columns B and C) -------------------------------
$bestData = @(
    ,@(2, 'B', 'After 1 laps have occured')
    ,@(2, 'C', ' the best unit found is: ')
    ,@(3, 'B', 'Killer has cost: 3 hp: 4 max Damage: 1 attack range: 3 move speed: 11 attack speed: 1')
    ,@(4, 'B', 'The total score was : 1.2708803')
    ,@(5, 'B', 'The score for round 1 was: 0.27088037')
    ,@(6, 'B', 'The unit was made in: 1 games')
    ,@(7, 'B', 'The unit was made in and was won with in : 0 games')
    ,@(8, 'B', 'The score for round 2 was: 1.0')
    ,@(9, 'B', 'The unit was made by player 0 in: 1 games')
    ,@(10, 'B', 'The unit was made by player 1 in: 1 games')
    ,@(11, 'B', 'The unit was made by both players in the same game in: 0 games')
    ,@(12, 'B', 'The unit was made in and was won with by player 0 in : 1 games')
    ,@(13, 'B', 'The unit was made in and was won with by player 1 in : 1 games')
    ,@(14, 'B', 'After 2 laps have occured')
    ,@(14, 'C', ' the best unit found is: ')
    ,@(15, 'B', 'Killer has cost: 3 hp: 4 max Damage: 1 attack range: 2 move speed: 11 attack speed: 1')
    ,@(16, 'B', 'The total score was : 1.0364103')
    ,@(17, 'B', 'The score for round 1 was: 0.20307693')
    ,@(18, 'B', 'The unit was made in: 1 games')
    ,@(19, 'B', 'The unit was made in and was won with in : 0 games')
    ,@(20, 'B', 'The score for round 2 was: 0.8333334')
    ,@(21, 'B', 'The unit was made by player 0 in: 1 games')
    ,@(22, 'B', 'The unit was made by player 1 in: 2 games')
    ,@(23, 'B', 'The unit was made by both players in the same game in: 0 games')
    ,@(24, 'B', 'The unit was made in and was won with by player 0 in : 1 games')
    ,@(25, 'B', 'The unit was made in and was won with by player 1 in : 2 games')
    ,@(26, 'B', 'After 3 laps have occured')
    ,@(26, 'C', ' the best unit found is: ')
    ,@(27, 'B', 'Killer has cost: 3 hp: 4 max Damage: 1 attack range: 2 move speed: 6 attack speed: 1')
    ,@(28, 'B', 'The total score was : 1.3181102')
    ,@(29, 'B', 'The score for round 1 was: 0.31811023')
    ,@(30, 'B', 'The unit was made in: 1 games')
    ,@(31, 'B', 'The unit was made in and was won with in : 0 games')
    ,@(32, 'B', 'The score for round 2 was: 1.0')
    ,@(33, 'B', 'The unit was made by player 0 in: 2 games')
    ,@(34, 'B', 'The unit was made by player 1 in: 1 games')
    ,@(35, 'B', 'The unit was made by both players in the same game in: 1 games')
    ,@(36, 'B', 'The unit was made in and was won with by player 0 in : 1 games')
    ,@(37, 'B', 'The unit was made in and was won with by player 1 in : 1 games')
    ,@(38, 'B', 'After 4 laps have occured')
    ,@(38, 'C', ' the best unit found is: ')
    ,@(39, 'B', 'Killer has cost: 3 hp: 4 max Damage: 1 attack range: 2 move speed: 6 attack speed: 1')
    ,@(40, 'B', 'The total score was : 1.3181102')
    ,@(41, 'B', 'The score for round 1 was: 0.31811023')
    ,@(42, 'B', 'The unit was made in: 1 games')
    ,@(43, 'B', 'The unit was made in and was won with in : 0 games')
    ,@(44, 'B', 'The score for round 2 was: 1.0')
    ,@(45, 'B', 'The unit was made by player 0 in: 2 games')
    ,@(46, 'B', 'The unit was made by player 1 in: 1 games')
    ,@(47, 'B', 'The unit was made by both players in the same game in: 1 games')
    ,@(48, 'B', 'The unit was made in and was won with by player 0 in : 1 games')
    ,@(49, 'B', 'The unit was made in and was won with by player 1 in : 1 games')
    ,@(50, 'B', 'After 5 laps have occured')
    ,@(50, 'C', ' the best unit found is: ')
    ,@(51, 'B', 'Killer has cost: 3 hp: 4 max Damage: 1 attack range: 2 move speed: 6 attack speed: 1')
    ,@(52, 'B', 'The total score was : 1.3181102')
    ,@(53, 'B', 'The score for round 1 was: 0.31811023')
    ,@(54, 'B', 'The unit was made in: 1 games')
    ,@(55, 'B', 'The unit was made in and was won with in : 0 games')
    ,@(56, 'B', 'The score for round 2 was: 1.0')
    ,@(57, 'B', 'The unit was made by player 0 in: 2 games')
    ,@(58, 'B', 'The unit was made by player 1 in: 1 games')
    ,@(59, 'B', 'The unit was made by both players in the same game in: 1 games')
    ,@(60, 'B', 'The unit was made in and was won with by player 0 in : 1 games')
    ,@(61, 'B', 'The unit was made in and was won with by player 1 in : 1 games')
    ,@(62, 'B', 'After 6 laps have occured')
    ,@(62, 'C', ' the best unit found is: ')
    ,@(63, 'B', 'Killer has cost: 3 hp: 4 max Damage: 1 attack range: 2 move speed: 6 attack speed: 1')
    ,@(64, 'B', 'The total score was : 1.3181102')
    ,@(65, 'B', 'The score for round 1 was: 0.31811023')
    ,@(66, 'B', 'The unit was made in: 1 games')
    ,@(67, 'B', 'The unit was made in and was won with in : 0 games')
    ,@(68, 'B', 'The score for round 2 was: 1.0')
    ,@(69, 'B', 'The unit was made by player 0 in: 2 games')
    ,@(70, 'B', 'The unit was made by player 1 in: 1 games')
    ,@(71, 'B', 'The unit was made by both players in the same game in: 1 games')
    ,@(72, 'B', 'The unit was made in and was won with by player 0 in : 1 games')
    ,@(73, 'B', 'The unit was made in and was won with by player 1 in : 1 games')
    ,@(74, 'B', 'After 1 laps have occured')
    ,@(74, 'C', ' the best unit found is: ')
    ,@(75, 'B', 'Killer has cost: 1 hp: 4 max Damage: 2 attack range: 1 move speed: 14 attack speed: 1')
    ,@(76, 'B', 'The total score was : 0.464634')
    ,@(77, 'B', 'The score for round 1 was: -0.535366')
    ,@(78, 'B', 'The unit was made in: 3 games')
    ,@(79, 'B', 'The unit was made in and was won with in : 3 games')
    ,@(80, 'B', 'The score for round 2 was: 1.0')
    ,@(81, 'B', 'The unit was made by player 0 in: 7 games')
    ,@(82, 'B', 'The unit was made by player 1 in: 5 games')
    ,@(83, 'B', 'The unit was made by both players in the same game in: 4 games')
    ,@(84, 'B', 'The unit was made in and was won with by player 0 in : 4 games')
    ,@(85, 'B', 'The unit was made in and was won with by player 1 in : 4 games')
    ,@(86, 'B', 'After 2 laps have occured')
    ,@(86, 'C', ' the best unit found is: ')
    ,@(87, 'B', 'Killer has cost: 1 hp: 4 max Damage: 2 attack range: 1 move speed: 14 attack speed: 1')
    ,@(88, 'B', 'The total score was : 1.2068771')
    ,@(89, 'B', 'The score for round 1 was: 0.20687713')
    ,@(90, 'B', 'The unit was made in: 6 games')
    ,@(91, 'B', 'The unit was made in and was won with in : 2 games')
    ,@(92, 'B', 'The score for round 2 was: 1.0')
    ,@(93, 'B', 'The unit was made by player 0 in: 9 games')
    ,@(94, 'B', 'The unit was made by player 1 in: 9 games')
    ,@(95, 'B', 'The unit was made by both players in the same game in: 8 games')
    ,@(96, 'B', 'The unit was made in and was won with by player 0 in : 5 games')
    ,@(97, 'B', 'The unit was made in and was won with by player 1 in : 5 games')
    ,@(98, 'B', 'After 3 laps have occured')
    ,@(98, 'C', ' the best unit found is: ')
    ,@(99, 'B', 'Killer has cost: 1 hp: 4 max Damage: 2 attack range: 1 move speed: 14 attack speed: 1')
    ,@(100, 'B', 'The total score was : 0.96403915')
    ,@(101, 'B', 'The score for round 1 was: 0.03546777')
    ,@(102, 'B', 'The unit was made in: 6 games')
    ,@(103, 'B', 'The unit was made in and was won with in : 2 games')
    ,@(104, 'B', 'The score for round 2 was: 0.9285714')
    ,@(105, 'B', 'The unit was made by player 0 in: 6 games')
    ,@(106, 'B', 'The unit was made by player 1 in: 6 games')
    ,@(107, 'B', 'The unit was made by both players in the same game in: 5 games')
    ,@(108, 'B', 'The unit was made in and was won with by player 0 in : 4 games')
    ,@(109, 'B', 'The unit was made in and was won with by player 1 in : 3 games')
)

foreach ($item in $bestData) {
    $row = $item[0]
    $col = $item[1]
    $val = $item[2]
    $bestResults5.Range("$col$row").Value = $val
}

# --- Selections: select full used range on RawResults#5 first (temporarily
#     activates it), then finally select on BestResults#5 so it becomes the
#     active/selected sheet in the saved workbook. ---------------------------
$rawResults5.Range("A1:A48").Select() | Out-Null

$bestResults5.Range("B1:C109").Select() | Out-Null
